# Update the "Metadata" sheet (sheet 1) of the FHIR StructureDefinition
# export workbook: bump the version/date, fill in the Publisher, swap the
# (duplicated) "Contact" row for a single "Jurisdiction" row, and refresh
# the root Extension's Short/Definition text on the "Elements" sheet
# (sheet 2) to match.

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item(1)
$wsElements = $wb.Worksheets.Item(2)

# Row 11 duplicated row 10's "Contact" / "No display for ContactDetail" —
# remove it so everything below shifts up by one row.
$wsMeta.Rows(11).Delete()

# Version bump.
$wsMeta.Range("B3").Value2 = "6.0.0"

# Refreshed publication date.
$wsMeta.Range("B8").Value2 = "2022-01-21T20:46:54+00:00"

# Publisher now has a value.
$wsMeta.Range("B9").Value2 = "Alvearie Team"

# Former "Contact" row (row 10, after the delete above) becomes "Jurisdiction".
$wsMeta.Range("A10").Value2 = "Jurisdiction"
$wsMeta.Range("B10").Value2 = "United States of America"

# The root Extension element's Short/Definition columns (K/L) on the
# Elements sheet now carry the profile-specific text instead of the
# generic "Extension" / "An Extension" placeholders.
$wsElements.Range("K2").Value2 = "Provider Region"
$wsElements.Range("L2").Value2 = "Customer-specific description of the region for the address"
